$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update median_rte (B) and count (C) values for 2024-08 (row 5), 2024-09 (row 6),
# and 2025-05 (row 14) as part of adding data for May and June.

$ws.Range("B5").Value = 0.8715745681574938
$ws.Range("C5").Value = 750

$ws.Range("B6").Value = 0.8543972907186217
$ws.Range("C6").Value = 714

$ws.Range("B14").Value = 0.8257090899136454
$ws.Range("C14").Value = 672
